$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 1 and push all existing data (and blank rows)
# down by one, keeping formulas/styling intact.
$ws.Rows.Item(1).Insert()

# Put the "how to use" instructions in the freshly inserted A1.
$ws.Range("A1").Value = "How to use: Paste file names into column B, then copy column D into code"

# Restore the view to the top of the sheet with A2 selected/active.
[void]$ws.Range("A2").Select()
